$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 already has the exact shape we need for the new row 21 (a plain
# numeric phone, a "2020-08-16" birthday, 0 points). Copying it and
# inserting the copy at row 21 both supplies those values/formats in one
# shot AND pushes the existing row 21 (phone "09876543", birthday
# 2020-08-16, points 0) down to row 22.
$ws.Rows.Item(15).Copy()
$ws.Rows.Item(21).Insert()

# Row 22 (the old row 21) already carries the right phone/points; only the
# birthday needs to be wiped out.
$ws.Cells.Item(22, 2).Value = ""
